$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 800
$ws.Range("J32").Value = 800
$ws.Range("L32").Value = 800
$ws.Range("N32").Value = -1452
$ws.Range("H55").Value = 46.22222
$ws.Range("I55").Value = 50.25
$ws.Range("J55").Value = 43
$ws.Range("K55").Value = 50.25
$ws.Range("L55").Value = 43
$ws.Range("M55").Value = 163.75
$ws.Range("N55").Value = -471
$ws.Range("H80").Value = 1431.5454
$ws.Range("J80").Value = 1537.125
$ws.Range("L80").Value = 4611.375
$ws.Range("N80").Value = -6607.375
$ws.Range("H83").Value = 1431.5454
$ws.Range("J83").Value = 1537.125
$ws.Range("L83").Value = 13834.125
$ws.Range("N83").Value = -23818.125
$ws.Range("H98").Value = 1716.3334
$ws.Range("I98").Value = 1716.3334
$ws.Range("K98").Value = 1716.3334
$ws.Range("M98").Value = -218.3334
$ws.Range("H107").Value = 227.45454
$ws.Range("I107").Value = 105.52631
$ws.Range("K107").Value = 105.52631
$ws.Range("M107").Value = 1814.47369
$ws.Range("H113").Value = 1493
$ws.Range("I113").Value = 1528.9166
$ws.Range("J113").Value = 1349.3334
$ws.Range("K113").Value = 1528.9166
$ws.Range("L113").Value = 1349.3334
$ws.Range("M113").Value = 1725.0834
$ws.Range("N113").Value = -7857.3334
$ws.Range("H116").Value = 5099.6
$ws.Range("I116").Value = 4374.5
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 4374.5
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = -932.5
$ws.Range("N116").Value = -14884
$ws.Range("H122").Value = 1716.3334
$ws.Range("I122").Value = 1716.3334
$ws.Range("K122").Value = 5149.0002
$ws.Range("M122").Value = -2699.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7711.795
$ws.Range("I32").Value = 6409.6113
$ws.Range("K32").Value = 6409.6113
$ws.Range("M32").Value = -6122.6113
$ws.Range("H61").Value = 2420.2666
$ws.Range("I61").Value = 2420.2666
$ws.Range("K61").Value = 2420.2666
$ws.Range("M61").Value = -2208.2666
$ws.Range("H136").Value = 2420.2666
$ws.Range("I136").Value = 2420.2666
$ws.Range("K136").Value = 7260.7998
$ws.Range("M136").Value = -4710.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2641.0833
$ws.Range("I22").Value = 1339
$ws.Range("J22").Value = 3943.1667
$ws.Range("K22").Value = 1339
$ws.Range("L22").Value = 3943.1667
$ws.Range("M22").Value = -989
$ws.Range("N22").Value = -4643.1667
$ws.Range("H31").Value = 5440.1055
$ws.Range("I31").Value = 3698.923
$ws.Range("K31").Value = 3698.923
$ws.Range("M31").Value = -3403.923
$ws.Range("H34").Value = 5440.1055
$ws.Range("I34").Value = 3698.923
$ws.Range("K34").Value = 3698.923
$ws.Range("M34").Value = -3496.923
$ws.Range("H88").Value = 13666.333
$ws.Range("J88").Value = 13666.333
$ws.Range("L88").Value = 13666.333
$ws.Range("N88").Value = -14478.333
$ws.Range("H91").Value = 13666.333
$ws.Range("J91").Value = 13666.333
$ws.Range("L91").Value = 13666.333
$ws.Range("N91").Value = -16474.333
$ws.Range("H132").Value = 2402.4707
$ws.Range("I132").Value = 1927.9375
$ws.Range("K132").Value = 5783.8125
$ws.Range("M132").Value = -3253.8125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 208.375
$ws.Range("I12").Value = 118
$ws.Range("J12").Value = 249.45454
$ws.Range("K12").Value = 354
$ws.Range("L12").Value = 748.3636200000001
$ws.Range("M12").Value = -181
$ws.Range("N12").Value = -1094.36362
$ws.Range("H18").Value = 1600
$ws.Range("I18").Value = 466.66666
$ws.Range("K18").Value = 1399.99998
$ws.Range("M18").Value = -1230.99998
$ws.Range("H55").Value = 3880
$ws.Range("J55").Value = 5916.6665
$ws.Range("L55").Value = 17749.9995
$ws.Range("N55").Value = -18103.9995
$ws.Range("H86").Value = 418.75
$ws.Range("J86").Value = 600
$ws.Range("L86").Value = 1800
$ws.Range("N86").Value = -4172
$ws.Range("H89").Value = 418.75
$ws.Range("J89").Value = 600
$ws.Range("L89").Value = 5400
$ws.Range("N89").Value = -17256
$ws.Range("H131").Value = 1000
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2040
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3097.875
$ws.Range("I122").Value = 2398.2
$ws.Range("J122").Value = 4264
$ws.Range("K122").Value = 7194.599999999999
$ws.Range("L122").Value = 12792
$ws.Range("M122").Value = -4744.599999999999
$ws.Range("N122").Value = -17692
$ws.Range("H126").Value = 4585.5293
$ws.Range("I126").Value = 3586.5
$ws.Range("K126").Value = 10759.5
$ws.Range("M126").Value = -8289.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1002.0769
$ws.Range("I22").Value = 743.6667
$ws.Range("J22").Value = 1223.5714
$ws.Range("K22").Value = 743.6667
$ws.Range("L22").Value = 1223.5714
$ws.Range("M22").Value = -448.6667
$ws.Range("N22").Value = -1813.5714
$ws.Range("H27").Value = 1002.0769
$ws.Range("I27").Value = 743.6667
$ws.Range("J27").Value = 1223.5714
$ws.Range("K27").Value = 743.6667
$ws.Range("L27").Value = 1223.5714
$ws.Range("M27").Value = -636.6667
$ws.Range("N27").Value = -1437.5714
$ws.Range("H46").Value = 8066.875
$ws.Range("I46").Value = 9824.6
$ws.Range("K46").Value = 9824.6
$ws.Range("M46").Value = -9636.6
$ws.Range("H93").Value = 1592.0714
$ws.Range("J93").Value = 2168.6
$ws.Range("L93").Value = 2168.6
$ws.Range("N93").Value = -4664.6
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3557
$ws.Range("I81").Value = 3557
$ws.Range("K81").Value = 7114
$ws.Range("M81").Value = -6053
$ws.Range("H84").Value = 3557
$ws.Range("I84").Value = 3557
$ws.Range("K84").Value = 35570
$ws.Range("M84").Value = -30266
$ws.Range("H113").Value = 608.625
$ws.Range("I113").Value = 511.625
$ws.Range("J113").Value = 705.625
$ws.Range("K113").Value = 1534.875
$ws.Range("L113").Value = 2116.875
$ws.Range("M113").Value = 635.125
$ws.Range("N113").Value = -6456.875
